# Add table-total support: a new placeholder row under the existing
# ${name} / ${date} / ${amount} rows that the templating engine can
# substitute with the computed grand total.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "`${totalAmount}"

# Reflect the user's final view state: zoomed in on the sheet and the
# newly added cell selected/active.
$excel.ActiveWindow.Zoom = 385
$ws.Range("A4").Select() | Out-Null

$wb.Save() | Out-Null
